$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.062481119763447
$ws.Cells.Item(2, 3).Value = 0.4492183535268595
$ws.Cells.Item(2, 5).Value = 0.2928396493372105
$ws.Cells.Item(2, 6).Value = 1.682642749130778
$ws.Cells.Item(2, 7).Value = 0.4784949337082267
$ws.Cells.Item(2, 8).Value = 0.6361308971360771
$ws.Cells.Item(2, 10).Value = 0.0483418449920876
$ws.Cells.Item(2, 12).Value = 0.4800696194376144
$ws.Cells.Item(2, 13).Value = 0.317045211634202
$ws.Cells.Item(2, 15).Value = 2.189351420827151

$ws.Cells.Item(3, 2).Value = 0.9563542948799295
$ws.Cells.Item(3, 3).Value = 0.4431311142110559
$ws.Cells.Item(3, 5).Value = 0.2950125084929827
$ws.Cells.Item(3, 6).Value = 1.692471632030362
$ws.Cells.Item(3, 7).Value = 0.4837530612755643
$ws.Cells.Item(3, 8).Value = 0.6428333587855946
$ws.Cells.Item(3, 10).Value = 0.04549827627575809
$ws.Cells.Item(3, 12).Value = 0.473106036226099
$ws.Cells.Item(3, 13).Value = 0.2973944632485299
$ws.Cells.Item(3, 15).Value = 2.214259556163412

$ws.Cells.Item(4, 2).Value = 0.8910989728689174
$ws.Cells.Item(4, 3).Value = 0.4394622221171431
$ws.Cells.Item(4, 5).Value = 0.2964481363376272
$ws.Cells.Item(4, 6).Value = 1.699491885536233
$ws.Cells.Item(4, 7).Value = 0.4874030058840546
$ws.Cells.Item(4, 8).Value = 0.6472857173037667
$ws.Cells.Item(4, 10).Value = 0.04374678443372204
$ws.Cells.Item(4, 12).Value = 0.468970624168648
$ws.Cells.Item(4, 13).Value = 0.2853636662090011
$ws.Cells.Item(4, 15).Value = 2.231142910692668

$ws.Cells.Item(5, 2).Value = 0.8644853883645283
$ws.Cells.Item(5, 3).Value = 0.4379846175157098
$ws.Cells.Item(5, 5).Value = 0.2970587206583826
$ws.Cells.Item(5, 6).Value = 1.702600504107743
$ws.Cells.Item(5, 7).Value = 0.4889961791713588
$ws.Cells.Item(5, 8).Value = 0.6491848392811974
$ws.Cells.Item(5, 10).Value = 0.04303169755921488
$ws.Cells.Item(5, 12).Value = 0.4673208865409038
$ws.Cells.Item(5, 13).Value = 0.2804701517479558
$ws.Cells.Item(5, 15).Value = 2.238422363290837

$ws.Cells.Item(6, 2).Value = 0.8600649793139041
$ws.Cells.Item(6, 3).Value = 0.4377403259133388
$ws.Cells.Item(6, 5).Value = 0.2971616521792768
$ws.Cells.Item(6, 6).Value = 1.703131656370942
$ws.Cells.Item(6, 7).Value = 0.489267108571724
$ws.Cells.Item(6, 8).Value = 0.6495053058118643
$ws.Cells.Item(6, 10).Value = 0.04291287850789871
$ws.Cells.Item(6, 12).Value = 0.467049097676508
$ws.Cells.Item(6, 13).Value = 0.279658149992585
$ws.Cells.Item(6, 15).Value = 2.239655219220978

$ws.Cells.Item(7, 2).Value = 0.8907401373563175
$ws.Cells.Item(7, 3).Value = 0.4394422234730655
$ws.Cells.Item(7, 5).Value = 0.2964562673756337
$ws.Cells.Item(7, 6).Value = 1.699532806019349
$ws.Cells.Item(7, 7).Value = 0.4874240638733696
$ws.Cells.Item(7, 8).Value = 0.6473109863742437
$ws.Cells.Item(7, 10).Value = 0.04373714587906363
$ws.Cells.Item(7, 12).Value = 0.4689482313106765
$ws.Cells.Item(7, 13).Value = 0.2852976330879642
$ws.Cells.Item(7, 15).Value = 2.231239467501624

$ws.Cells.Item(8, 2).Value = 1.025908988990409
$ws.Cells.Item(8, 3).Value = 0.4471053467059249
$ws.Cells.Item(8, 5).Value = 0.2935678096361638
$ws.Cells.Item(8, 6).Value = 1.685827282490251
$ws.Cells.Item(8, 7).Value = 0.4802203585855125
$ws.Cells.Item(8, 8).Value = 0.6383719596975865
$ws.Cells.Item(8, 10).Value = 0.04736255950193424
$ws.Cells.Item(8, 12).Value = 0.4776395873242194
$ws.Cells.Item(8, 13).Value = 0.3102626182317749
$ws.Cells.Item(8, 15).Value = 2.19760958201914

$ws.Cells.Item(9, 2).Value = 1.29016719727224
$ws.Cells.Item(9, 3).Value = 0.4626689387729925
$ws.Cells.Item(9, 5).Value = 0.2887071927047025
$ws.Cells.Item(9, 6).Value = 1.6667673777456
$ws.Cells.Item(9, 7).Value = 0.4694461472930271
$ws.Cells.Item(9, 8).Value = 0.6235165026145992
$ws.Cells.Item(9, 10).Value = 0.05442623233900434
$ws.Cells.Item(9, 12).Value = 0.4957882420772535
$ws.Cells.Item(9, 13).Value = 0.3594815520707186
$ws.Cells.Item(9, 15).Value = 2.144292368163391

$ws.Cells.Item(10, 2).Value = 1.483751573140978
$ws.Cells.Item(10, 3).Value = 0.4744199739757562
$ws.Cells.Item(10, 5).Value = 0.2856238427921429
$ws.Cells.Item(10, 6).Value = 1.657529340279027
$ws.Cells.Item(10, 7).Value = 0.4635857379921262
$ws.Cells.Item(10, 8).Value = 0.6142320895939193
$ws.Cells.Item(10, 10).Value = 0.05958599788876739
$ws.Cells.Item(10, 12).Value = 0.5097867408118475
$ws.Cells.Item(10, 13).Value = 0.3957880678451886
$ws.Cells.Item(10, 15).Value = 2.112845277461673

$ws.Cells.Item(11, 2).Value = 1.571680893625512
$ws.Cells.Item(11, 3).Value = 0.4798324565115593
$ws.Cells.Item(11, 5).Value = 0.2843265961270571
$ws.Cells.Item(11, 6).Value = 1.65436157425664
$ws.Cells.Item(11, 7).Value = 0.4613685323534895
$ws.Cells.Item(11, 8).Value = 0.6103622492414118
$ws.Cells.Item(11, 10).Value = 0.06192641866002901
$ws.Cells.Item(11, 12).Value = 0.5162975781407795
$ws.Cells.Item(11, 13).Value = 0.4123335339681944
$ws.Cells.Item(11, 15).Value = 2.100221989153368

$ws.Cells.Item(12, 2).Value = 1.604956680262376
$ws.Cells.Item(12, 3).Value = 0.4818914408328396
$ws.Cells.Item(12, 5).Value = 0.2838504814583231
$ws.Cells.Item(12, 6).Value = 1.653310790517367
$ws.Cells.Item(12, 7).Value = 0.4605936612480335
$ws.Cells.Item(12, 8).Value = 0.6089477039992914
$ws.Cells.Item(12, 10).Value = 0.06281165590897331
$ws.Cells.Item(12, 12).Value = 0.5187834202519355
$ws.Cells.Item(12, 13).Value = 0.4186027861953647
$ws.Cells.Item(12, 15).Value = 2.09568420842362

$ws.Cells.Item(13, 2).Value = 1.597791115872099
$ws.Cells.Item(13, 3).Value = 0.4814475874006234
$ws.Cells.Item(13, 5).Value = 0.2839523490774241
$ws.Cells.Item(13, 6).Value = 1.653530478611756
$ws.Cells.Item(13, 7).Value = 0.4607576610942274
$ws.Cells.Item(13, 8).Value = 0.6092500884427849
$ws.Cells.Item(13, 10).Value = 0.06262105086355518
$ws.Cells.Item(13, 12).Value = 0.5182471488973306
$ws.Cells.Item(13, 13).Value = 0.4172524258363168
$ws.Cells.Item(13, 15).Value = 2.096650712762042

$ws.Cells.Item(14, 2).Value = 1.574418947603704
$ws.Cells.Item(14, 3).Value = 0.4800016633722919
$ws.Cells.Item(14, 5).Value = 0.2842871229392081
$ws.Cells.Item(14, 6).Value = 1.654272143886814
$ws.Cells.Item(14, 7).Value = 0.4613034846243806
$ws.Cells.Item(14, 8).Value = 0.6102448539434064
$ws.Cells.Item(14, 10).Value = 0.06199926859837035
$ws.Cells.Item(14, 12).Value = 0.5165016836512564
$ws.Cells.Item(14, 13).Value = 0.4128492343316026
$ws.Cells.Item(14, 15).Value = 2.099843802641686

$ws.Cells.Item(15, 2).Value = 1.56010000805793
$ws.Cells.Item(15, 3).Value = 0.479117210063265
$ws.Cells.Item(15, 5).Value = 0.2844941503711258
$ws.Cells.Item(15, 6).Value = 1.654745810776177
$ws.Cells.Item(15, 7).Value = 0.4616462541278139
$ws.Cells.Item(15, 8).Value = 0.6108608031073786
$ws.Cells.Item(15, 10).Value = 0.06161827336460846
$ws.Cells.Item(15, 12).Value = 0.5154351769766095
$ws.Cells.Item(15, 13).Value = 0.4101526413226466
$ws.Cells.Item(15, 15).Value = 2.101831242594727

$ws.Cells.Item(16, 2).Value = 1.478002335085932
$ws.Cells.Item(16, 3).Value = 0.4740675822013429
$ws.Cells.Item(16, 5).Value = 0.2857107389809155
$ws.Cells.Item(16, 6).Value = 1.657757181117887
$ws.Cells.Item(16, 7).Value = 0.4637396851086564
$ws.Cells.Item(16, 8).Value = 0.614492111513151
$ws.Cells.Item(16, 10).Value = 0.05943290478403895
$ws.Cells.Item(16, 12).Value = 0.5093640986200398
$ws.Cells.Item(16, 13).Value = 0.3947073418057059
$ws.Cells.Item(16, 15).Value = 2.113704130479192

$ws.Cells.Item(17, 2).Value = 1.427602542145848
$ws.Cells.Item(17, 3).Value = 0.4709867690542922
$ws.Cells.Item(17, 5).Value = 0.2864840470714061
$ws.Cells.Item(17, 6).Value = 1.65986956006968
$ws.Cells.Item(17, 7).Value = 0.4651390165561935
$ws.Cells.Item(17, 8).Value = 0.6168104030842727
$ws.Cells.Item(17, 10).Value = 0.05809047627521835
$ws.Cells.Item(17, 12).Value = 0.5056761246878949
$ws.Cells.Item(17, 13).Value = 0.3852394143260085
$ws.Cells.Item(17, 15).Value = 2.121418954116237

$ws.Cells.Item(18, 2).Value = 1.398601496340916
$ws.Cells.Item(18, 3).Value = 0.4692210747800232
$ws.Cells.Item(18, 5).Value = 0.2869387545067834
$ws.Cells.Item(18, 6).Value = 1.661181932305269
$ws.Cells.Item(18, 7).Value = 0.4659860985322979
$ws.Cells.Item(18, 8).Value = 0.6181771124770705
$ws.Cells.Item(18, 10).Value = 0.05731771171337385
$ws.Cells.Item(18, 12).Value = 0.5035683579953485
$ws.Cells.Item(18, 13).Value = 0.3797965080907062
$ws.Cells.Item(18, 15).Value = 2.126014619384549

$ws.Cells.Item(19, 2).Value = 1.388780168526523
$ws.Cells.Item(19, 3).Value = 0.4686243312524425
$ws.Cells.Item(19, 5).Value = 0.2870944154917758
$ws.Cells.Item(19, 6).Value = 1.661643005486383
$ws.Cells.Item(19, 7).Value = 0.4662801514246127
$ws.Cells.Item(19, 8).Value = 0.6186455736786272
$ws.Cells.Item(19, 10).Value = 0.05705595946635356
$ws.Cells.Item(19, 12).Value = 0.5028570223797857
$ws.Cells.Item(19, 13).Value = 0.377954125614977
$ws.Cells.Item(19, 15).Value = 2.127597806240956

$ws.Cells.Item(20, 2).Value = 1.432968985791888
$ws.Cells.Item(20, 3).Value = 0.4713140756406347
$ws.Cells.Item(20, 5).Value = 0.286400700526583
$ws.Cells.Item(20, 6).Value = 1.659634614473049
$ws.Cells.Item(20, 7).Value = 0.4649856838257946
$ws.Cells.Item(20, 8).Value = 0.6165601715165252
$ws.Cells.Item(20, 10).Value = 0.05823344621387605
$ws.Cells.Item(20, 12).Value = 0.5060673244590959
$ws.Cells.Item(20, 13).Value = 0.3862470049170241
$ws.Cells.Item(20, 15).Value = 2.120581310639494

$ws.Cells.Item(21, 2).Value = 1.581284509439968
$ws.Cells.Item(21, 3).Value = 0.4804261130998952
$ws.Cells.Item(21, 5).Value = 0.2841883815029504
$ws.Cells.Item(21, 6).Value = 1.654050260990388
$ws.Cells.Item(21, 7).Value = 0.4611414044212765
$ws.Cells.Item(21, 8).Value = 0.6099512861956171
$ws.Cells.Item(21, 10).Value = 0.06218192948561097
$ws.Cells.Item(21, 12).Value = 0.5170138190614466
$ws.Cells.Item(21, 13).Value = 0.4141424572151564
$ws.Cells.Item(21, 15).Value = 2.098899332582263

$ws.Cells.Item(22, 2).Value = 1.678093152988424
$ws.Cells.Item(22, 3).Value = 0.4864360142370003
$ws.Cells.Item(22, 5).Value = 0.2828306465112096
$ws.Cells.Item(22, 6).Value = 1.651267767012854
$ws.Cells.Item(22, 7).Value = 0.4590063762857
$ws.Cells.Item(22, 8).Value = 0.605928570230148
$ws.Cells.Item(22, 10).Value = 0.06475646715535532
$ws.Cells.Item(22, 12).Value = 0.5242863828846538
$ws.Cells.Item(22, 13).Value = 0.4323959693981507
$ws.Cells.Item(22, 15).Value = 2.086141935585147

$ws.Cells.Item(23, 2).Value = 1.626436568417603
$ws.Cells.Item(23, 3).Value = 0.4832234884570994
$ws.Cells.Item(23, 5).Value = 0.2835472398929291
$ws.Cells.Item(23, 6).Value = 1.652673488686105
$ws.Cells.Item(23, 7).Value = 0.4601112769073694
$ws.Cells.Item(23, 8).Value = 0.6080484265463895
$ws.Cells.Item(23, 10).Value = 0.06338295694020957
$ws.Cells.Item(23, 12).Value = 0.5203941153547333
$ws.Cells.Item(23, 13).Value = 0.4226518219365047
$ws.Cells.Item(23, 15).Value = 2.092821344221647

$ws.Cells.Item(24, 2).Value = 1.430542896003772
$ws.Cells.Item(24, 3).Value = 0.4711660831648032
$ws.Cells.Item(24, 5).Value = 0.2864383499673586
$ws.Cells.Item(24, 6).Value = 1.659740528310053
$ws.Cells.Item(24, 7).Value = 0.4650548729113595
$ws.Cells.Item(24, 8).Value = 0.6166731956356983
$ws.Cells.Item(24, 10).Value = 0.05816881256708228
$ws.Cells.Item(24, 12).Value = 0.5058904240956821
$ws.Cells.Item(24, 13).Value = 0.3857914721493501
$ws.Cells.Item(24, 15).Value = 2.120959509935162

$ws.Cells.Item(25, 2).Value = 1.218772489095727
$ws.Cells.Item(25, 3).Value = 0.4584022567019019
$ws.Cells.Item(25, 5).Value = 0.2899363084003586
$ws.Cells.Item(25, 6).Value = 1.671086753312252
$ws.Cells.Item(25, 7).Value = 0.4720007526404544
$ws.Cells.Item(25, 8).Value = 0.6272490552942571
$ws.Cells.Item(25, 10).Value = 0.05252043018045072
$ws.Cells.Item(25, 12).Value = 0.4907611840859403
$ws.Cells.Item(25, 13).Value = 0.3461399411743642
$ws.Cells.Item(25, 15).Value = 2.157361273123556
